# d20_JN.xlsx — series of logistic regressions for response at each split
#
# The pctcorrect column (I) for the "ALL" / musicians-unfiltered split
# (rows 3-6) had been computed from the wrong regression output; it's
# corrected here to match the pctrelated column (H), same as the other
# splits on the sheet already show.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 69.827586206896598
$ws.Range("I4").Value = 63.557993730407503
$ws.Range("I5").Value = 56.198347107437996
$ws.Range("I6").Value = 50.236034618410699

# Update the saved cursor/selection position for the sheet.
$ws.Range("K6").Select() | Out-Null
